$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") rows 2-43 all changed from 45799 (2025-05-22) to 45800 (2025-05-23)
$ws.Range("C2:C43").Value = 45800
